$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200, shifting existing rows 200-269 down to 201-270
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with data
$ws.Cells.Item(200, 1).Value = 8
$ws.Cells.Item(200, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(200, 3).Value = "Coquimbo"
$ws.Cells.Item(200, 4).Value = 44795
$ws.Cells.Item(200, 5).Value = 4
$ws.Cells.Item(200, 6).Value = 100112021
$ws.Cells.Item(200, 7).Value = "Ají"
$ws.Cells.Item(200, 8).Value = "Inferno"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 600
$ws.Cells.Item(200, 11).Value = 14500
$ws.Cells.Item(200, 12).Value = 15000
$ws.Cells.Item(200, 13).Value = 14750
$ws.Cells.Item(200, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(200, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(200, 16).Value = 1475
$ws.Cells.Item(200, 17).Value = 10
$ws.Cells.Item(200, 18).Value = "Hortaliza"

# Copy the style (date format) from D199 into the new D200 cell
$ws.Cells.Item(199, 4).Copy()
$ws.Cells.Item(200, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
